$d = $word.ActiveDocument

# The author removed one of the red "pending items" notes (and the related
# blank/italic paragraphs that went with it) from the section right after
# "...unirse a nosotros en esta misión." and right before the
# "Ah, y por cierto..." paragraph. Concretely this deletes:
#   - an empty spacer paragraph
#   - the red "***Valorar si incluir en artículo o reservarlo para nuestra
#     revisión más adelante:" paragraph
#   - the italic "No obstante... chatbot de OpenAI..." paragraph
#   - the empty (shaded) spacer paragraph that followed it
#
# Locate the paragraphs robustly by scanning paragraph text instead of
# relying on fixed paragraph indices.

$paras = $d.Paragraphs
$count = $paras.Count

$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*Valorar si incluir en artículo*") {
        $startIndex = $i - 1   # include the blank spacer paragraph right before it
    }
    if ($t -like "*Ah, y por cierto*") {
        $endIndex = $i - 1     # stop right before this paragraph
        break
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $pStart = $paras.Item($startIndex)
    $pEnd = $paras.Item($endIndex)
    $deleteRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $deleteRange.Delete()
}
